$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 209 ("Species name" / QUALITY / 1). This shifts all subsequent
# rows (210..230) up by one, so the sheet's last row becomes 229 instead of 230.
$ws.Rows.Item(209).Delete()
